$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 81; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value2()
    if ($current -eq 45184) {
        $cell.Value = 45185
    }
}
